$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.264.50'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.909.53'
$ws.Range("E3").Value = '  +0.30%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '307.92'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("E6").Value = '  +0.14%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5247'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.40%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3825'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +1.72%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.07314'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.92%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '21.63'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.43%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.9072'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.55%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08173'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -4.11%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '96.30'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +1.25%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '5.374'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").Value = '1.756.85'
$ws.Range("E15").Value = '  -7.55%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.000008689'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.71%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '14.77'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '27.294.69'
$ws.Range("E20").Value = '  +0.41%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '5.131'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +1.28%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '10.81'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +1.99%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.503'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +1.23%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '150.24'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +2.26%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '18.26'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.28%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '1.741'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.69%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '117.09'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.85%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '4.859'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.00%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.882'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.35%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.09232'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.21%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.8255'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.64%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.05085'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.63%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.988'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '3.374'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.747'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +4.88%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.5754'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.84%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.02005'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +0.47%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '1.083'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +0.73%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '9.059'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +0.04%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '6.621'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '117.24'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.84%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.1523'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.40%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.4942'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '10.16'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("E48").Value = '  +1.77%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '38.72'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.30%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.05988'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
